$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells: "_old" suffix -> "_FV2404", "_new" suffix -> "_FV2410" ---
$headers = @(
  "Segmentname",
  "Segmentgruppe",
  "Segment",
  "Datenelement",
  "Segment ID",
  "Code",
  "Qualifier",
  "Beschreibung",
  "Bedingungsausdruck",
  "Bedingung"
)

for ($i = 0; $i -lt $headers.Count; $i++) {
  $col = $i + 1
  $ws.Cells.Item(1, $col).Value2 = "$($headers[$i])_FV2404"
}

# column 11 ("diff") stays as-is

for ($i = 0; $i -lt $headers.Count; $i++) {
  $col = $i + 12
  $ws.Cells.Item(1, $col).Value2 = "$($headers[$i])_FV2410"
}

# --- 2. Freeze the header row (pane split after row 1) ---
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the used range into an Excel Table (ListObject) ---
$rng = $ws.Range("A1:U68")
$tbl = $ws.ListObjects.Add(1, $rng, $false, 1)
$tbl.Name = "Table1"
